# Apply the edits described by the diff to NIT-9017490953.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Copy "last row" formatting (borders) from row 45 onto row 29 before we
#        delete the intervening rows, so the new last data row keeps the
#        closing-border look that row 45 used to have. ---
$ws.Range("B45:J45").Copy()
$ws.Range("B29:J29").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- 2. Header area updates ---
$ws.Range("D2").Value = "ESTADO DE CUENTA"
$ws.Range("B7").Value = "RAZON SOCIAL:"
$ws.Range("B11").Value = "VALOR MORA"
$ws.Range("E11").Value = 787680
$ws.Range("B13").Value = "Cant. Trabajadores"
$ws.Range("C13").Value = 3
$ws.Range("E13").Value = "Cant. Periodos"

# --- 3. Column header row (row 15) - Novedad de Ingreso / Novedad de Retiro swap ---
$ws.Range("H15").Value = "Novedad de Ingreso"
$ws.Range("I15").Value = "Novedad de Retiro"
$ws.Range("J15").Value = "Observaciones"

# --- 4. Rewrite the data table (rows 16-29) with the new account data ---
$data = @(
    @("CC", "45514702",   "ERNESTINA PARRA LUNA",             "2403", 52000, 5800000),
    @("CC", "8865467",    "GUILLERMO ALONSO ROJAS MARTINEZ",  "2403", 52000, 1300000),
    @("CC", "1002196724", "MARCOS ANTONIO AYOLA CARO",        "2403", 57600, 1440000),
    @("CC", "1002196724", "MARCOS ANTONIO AYOLA CARO",        "2404", 57600, 1440000),
    @("CC", "45514702",   "ERNESTINA PARRA LUNA",             "2405", 52000, 5800000),
    @("CC", "1002196724", "MARCOS ANTONIO AYOLA CARO",        "2405", 57600, 1440000),
    @("CC", "1002196724", "MARCOS ANTONIO AYOLA CARO",        "2406", 57600, 1440000),
    @("CC", "1002196724", "MARCOS ANTONIO AYOLA CARO",        "2407", 57600, 1440000),
    @("CC", "1002196724", "MARCOS ANTONIO AYOLA CARO",        "2408", 57600, 1440000),
    @("CC", "1002196724", "MARCOS ANTONIO AYOLA CARO",        "2409", 57600, 1440000),
    @("CC", "1002196724", "MARCOS ANTONIO AYOLA CARO",        "2410", 57600, 1440000),
    @("CC", "1002196724", "MARCOS ANTONIO AYOLA CARO",        "2411", 57600, 1440000),
    @("CC", "1002196724", "MARCOS ANTONIO AYOLA CARO",        "2412", 57600, 1440000),
    @("CC", "1002196724", "MARCOS ANTONIO AYOLA CARO",        "2501", 55680, 1440000)
)

$r = 16
foreach ($row in $data) {
    $ws.Range("B$r").Value = $row[0]
    $ws.Range("C$r").Value = $row[1]
    $ws.Range("D$r").Value = $row[2]
    $ws.Range("E$r").Value = $row[3]
    $ws.Range("F$r").Value = $row[4]
    $ws.Range("G$r").Value = $row[5]
    $r = $r + 1
}

# --- 5. Remove the now-obsolete extra data rows (old rows 30-45) ---
$ws.Rows("30:45").Delete()

# --- 6. Column D width ---
$ws.Columns("D").ColumnWidth = 36
